$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 83, shifting rows 83..148 down to 84..149.
$ws.Rows.Item(83).Insert()

# Populate the newly inserted row 83 with the new data record.
$ws.Cells.Item(83, 1).Value = 8
$ws.Cells.Item(83, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(83, 3).Value = "Coquimbo"
$ws.Cells.Item(83, 4).NumberFormat = $ws.Cells.Item(84, 4).NumberFormat
$ws.Cells.Item(83, 4).Value = 45040
$ws.Cells.Item(83, 5).Value = 4
$ws.Cells.Item(83, 6).Value = "Fruta"
$ws.Cells.Item(83, 7).Value = 100109
$ws.Cells.Item(83, 8).Value = "Uva"
$ws.Cells.Item(83, 9).Value = 100109001
$ws.Cells.Item(83, 10).Value = "Uva"
$ws.Cells.Item(83, 11).Value = "Red Globe"
$ws.Cells.Item(83, 12).Value = "Primera"
$ws.Cells.Item(83, 13).Value = 400
$ws.Cells.Item(83, 14).Value = 8000
$ws.Cells.Item(83, 15).Value = 9000
$ws.Cells.Item(83, 16).Value = 8500
$ws.Cells.Item(83, 17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(83, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(83, 19).Value = 472
$ws.Cells.Item(83, 20).Value = 18
